# Update ARCtrl MSEval_Cold assay test object:
# - Rename "Source Name" column header to "Input [Raw Data File]"
# - Rename "Sample Name" column header to "Output [Derived Data File]"
# These header cells drive the names of the corresponding columns in the
# "annotationTableNiceZebra52" Excel table on the MSEval worksheet, so
# updating the cell values also updates the table column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSEval")

$ws.Range("A1").Value = "Input [Raw Data File]"
$ws.Range("K1").Value = "Output [Derived Data File]"

# Match the author's final selection on the MSEval sheet.
$ws.Range("E7").Select()
